# chore: update Sheets via scheduled runner
#
# Refreshes the market-board derived columns (H: currentAveragePrice,
# I: currentAveragePriceNQ, J: currentAveragePriceHQ, K: LevePriceNQ,
# L: LevePriceHQ, M: LeveProfitNQ, N: LeveProfitHQ) for a handful of
# leve rows across several class sheets, as pulled by the scheduled
# price-refresh job.

$wb = $excel.ActiveWorkbook

function Set-LeveRow {
    param(
        [object]$ws,
        [int]$row,
        [object]$H = $null,
        [object]$I = $null,
        [object]$J = $null,
        [object]$K = $null,
        [object]$L = $null,
        [object]$M = $null,
        [object]$N = $null
    )

    if ($null -ne $H) { $ws.Range("H$row").Value = $H }
    if ($null -ne $I) { $ws.Range("I$row").Value = $I }
    if ($null -ne $J) { $ws.Range("J$row").Value = $J }
    if ($null -ne $K) { $ws.Range("K$row").Value = $K }
    if ($null -ne $L) { $ws.Range("L$row").Value = $L }
    if ($null -ne $M) { $ws.Range("M$row").Value = $M }
    if ($null -ne $N) { $ws.Range("N$row").Value = $N }
}

# --- ALC ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

Set-LeveRow -ws $ws -row 46 `
    -H 1775 -I 1166.6666 -J 1977.7778 -K 3499.9998 -L 5933.3334 `
    -M -3380.9998 -N -6171.3334

Set-LeveRow -ws $ws -row 60 `
    -H 1775 -I 1166.6666 -J 1977.7778 -K 3499.9998 -L 5933.3334 `
    -M -3015.9998 -N -6901.3334

Set-LeveRow -ws $ws -row 64 `
    -H 58827816 -I 142862480 -J 3550 -K 142862480 -L 3550 `
    -M -142862232 -N -4046

Set-LeveRow -ws $ws -row 67 `
    -H 58827816 -I 142862480 -J 3550 -K 142862480 -L 3550 `
    -M -142861622 -N -5266

# --- ARM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

Set-LeveRow -ws $ws -row 61 `
    -H 228903.95 -I 1588.1111 -J 1251825.2 -K 1588.1111 -L 1251825.2 `
    -M -1376.1111 -N -1252249.2

Set-LeveRow -ws $ws -row 136 `
    -H 228903.95 -I 1588.1111 -J 1251825.2 -K 4764.3333 -L 3755475.6 `
    -M -2214.3333 -N -3760575.6

# --- BSM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

Set-LeveRow -ws $ws -row 80 `
    -H 664.5217 -I 575.375 -J 712.06665 -K 575.375 -L 712.06665 `
    -M 422.625 -N -2708.06665

Set-LeveRow -ws $ws -row 83 `
    -H 664.5217 -I 575.375 -J 712.06665 -K 2876.875 -L 3560.33325 `
    -M 2115.125 -N -13544.33325

# --- CRP -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

Set-LeveRow -ws $ws -row 31 `
    -H 1664.0103 -I 1140.431 -J 2423.2 -K 1140.431 -L 2423.2 `
    -M -845.431 -N -3013.2

Set-LeveRow -ws $ws -row 34 `
    -H 1664.0103 -I 1140.431 -J 2423.2 -K 1140.431 -L 2423.2 `
    -M -938.431 -N -2827.2

Set-LeveRow -ws $ws -row 58 `
    -H 76924310 -I 125000700 -J 2091.2 -K 125000700 -L 2091.2 `
    -M -125000497 -N -2497.2

Set-LeveRow -ws $ws -row 93 `
    -H 5450.5 -I 4480.533 -J 20000 -K 4480.533 -L 20000 `
    -M -2608.533 -N -23744

Set-LeveRow -ws $ws -row 136 `
    -H 76924310 -I 125000700 -J 2091.2 -K 375002100 -L 6273.599999999999 `
    -M -374999550 -N -11373.6

# --- CUL -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

Set-LeveRow -ws $ws -row 3 `
    -H 4759.9688 -I 2853.913 -J 9631 -K 8561.739 -L 28893 `
    -M -8449.739 -N -29117

Set-LeveRow -ws $ws -row 139 `
    -H 368228.22 -I 550729.9399999999 -J 3224.8 -K 1652189.82 -L 9674.400000000001 `
    -M -1647049.82 -N -19954.4

Set-LeveRow -ws $ws -row 140 `
    -H 31943.676 -I 38645 -J 3223.7144 -K 115935 -L 9671.143199999999 `
    -M -110755 -N -20031.1432

Set-LeveRow -ws $ws -row 141 `
    -H 41751.56 -I 41751.56 -K 125254.68 -M -120074.68

# --- GSM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

Set-LeveRow -ws $ws -row 70 `
    -H 3965.5557 -I 3985.7144 -J 3895 -K 3985.7144 -L 3895 `
    -M -3715.7144 -N -4435

Set-LeveRow -ws $ws -row 73 `
    -H 3965.5557 -I 3985.7144 -J 3895 -K 3985.7144 -L 3895 `
    -M -3049.7144 -N -5767

Set-LeveRow -ws $ws -row 132 `
    -H 4466.7554 -J 2584.1875 -L 7752.5625 -N -12812.5625

# --- LTW -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

Set-LeveRow -ws $ws -row 22 `
    -H 405.45456 -I 294 -J 498.33334 -K 294 -L 498.33334 `
    -M 1 -N -1088.33334

Set-LeveRow -ws $ws -row 27 `
    -H 405.45456 -I 294 -J 498.33334 -K 294 -L 498.33334 `
    -M -187 -N -712.33334

Set-LeveRow -ws $ws -row 46 `
    -H 1495.3636 -I 1108.1666 -J 1960 -K 1108.1666 -L 1960 `
    -M -920.1666 -N -2336

Set-LeveRow -ws $ws -row 136 `
    -H 5579.2 -I 1647.6471 -J 10720.462 -K 4942.9413 -L 32161.386 `
    -M -2392.9413 -N -37261.386

# --- WVR -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 53 lost its HQ leve (the NQ/HQ market split collapsed to 0), so the
# LeveProfitHQ cell (N53) no longer applies and is cleared entirely rather
# than zeroed.
Set-LeveRow -ws $ws -row 53 -H 0 -J 0 -L 0
$ws.Range("N53").ClearContents()

Write-Output "Leve profit columns refreshed."
